$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section note above the new pattern blocks
$ws.Range("A20").Value = "medizinkram pattern erweitert"

# Row 21 - "hungry" pattern (tag first, then pattern, then response)
$ws.Range("B21").Value = "hungry"
$ws.Range("A21").Value = "i am hungry…"
$ws.Range("D21").Value = "i am also hungry"

# Row 19 - "like_simple" pattern block (tag first, then note, then pattern, then response)
$ws.Range("B19").Value = "like_simple"
$ws.Range("C19").Value = " falls wir like noch brauchen"
$ws.Range("A19").Value = "i like football "
$ws.Range("D19").Value = "That is very nice, tell me more"

# Row 22 - "dontLike" pattern (tag first, then pattern, then response)
$ws.Range("B22").Value = "dontLike"
$ws.Range("A22").Value = "I don't like football"
$ws.Range("D22").Value = "I agree with you"

# Match the workbook's post-edit selection
$ws.Range("D22").Select() | Out-Null
